$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1211.3269
$ws.Range("I15").Value = 1211.3269
$ws.Range("K15").Value = 3633.9807
$ws.Range("M15").Value = -3464.9807
$ws.Range("H28").Value = 1003.1667
$ws.Range("I28").Value = 614.6667
$ws.Range("J28").Value = 1391.6666
$ws.Range("K28").Value = 614.6667
$ws.Range("L28").Value = 1391.6666
$ws.Range("M28").Value = -129.6667
$ws.Range("N28").Value = -2361.6666
$ws.Range("H74").Value = 5750
$ws.Range("I74").Value = 5366.6665
$ws.Range("J74").Value = 6900
$ws.Range("K74").Value = 5366.6665
$ws.Range("L74").Value = 6900
$ws.Range("M74").Value = -4430.6665
$ws.Range("N74").Value = -8772
$ws.Range("H77").Value = 5750
$ws.Range("I77").Value = 5366.6665
$ws.Range("J77").Value = 6900
$ws.Range("K77").Value = 26833.3325
$ws.Range("L77").Value = 34500
$ws.Range("M77").Value = -22153.3325
$ws.Range("N77").Value = -43860
$ws.Range("H92").Value = 495.1579
$ws.Range("I92").Value = 576.9375
$ws.Range("K92").Value = 576.9375
$ws.Range("M92").Value = 671.0625
$ws.Range("H111").Value = 14024.728
$ws.Range("I111").Value = 22097.666
$ws.Range("J111").Value = 4337.2
$ws.Range("K111").Value = 66292.99800000001
$ws.Range("L111").Value = 13011.6
$ws.Range("M111").Value = -63225.99800000001
$ws.Range("N111").Value = -19145.6
$ws.Range("H129").Value = 801.8
$ws.Range("I129").Value = 484
$ws.Range("K129").Value = 1452
$ws.Range("M129").Value = 3548

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 1435000.1
$ws.Range("I6").Value = 5004001
$ws.Range("J6").Value = 7399.8
$ws.Range("K6").Value = 5004001
$ws.Range("L6").Value = 7399.8
$ws.Range("M6").Value = -5003828
$ws.Range("N6").Value = -7745.8
$ws.Range("H26").Value = 3508.3
$ws.Range("I26").Value = 2863.8333
$ws.Range("J26").Value = 4475
$ws.Range("K26").Value = 2863.8333
$ws.Range("L26").Value = 4475
$ws.Range("M26").Value = -2533.8333
$ws.Range("N26").Value = -5135
$ws.Range("H39").Value = 25000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 25000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 25000
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("H42").Value = 10750
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 10750
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 10750
$ws.Range("N42").Value = -11722
$ws.Range("H88").Value = 2567
$ws.Range("I88").Value = 2689.3333
$ws.Range("J88").Value = 2200
$ws.Range("K88").Value = 2689.3333
$ws.Range("L88").Value = 2200
$ws.Range("M88").Value = -2283.3333
$ws.Range("N88").Value = -3012
$ws.Range("H91").Value = 2567
$ws.Range("I91").Value = 2689.3333
$ws.Range("J91").Value = 2200
$ws.Range("K91").Value = 2689.3333
$ws.Range("L91").Value = 2200
$ws.Range("M91").Value = -1285.3333
$ws.Range("N91").Value = -5008
$ws.Range("N39").Value = -26040
$ws.Range("M39").Value = ""
$ws.Range("N40").Value = ""
$ws.Range("M42").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1326.6666
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 1326.6666
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 1326.6666
$ws.Range("H86").Value = 220950
$ws.Range("J86").Value = 1475
$ws.Range("L86").Value = 1475
$ws.Range("N86").Value = -3721
$ws.Range("H89").Value = 220950
$ws.Range("J89").Value = 1475
$ws.Range("L89").Value = 7375
$ws.Range("N89").Value = -18607
$ws.Range("N7").Value = -1552.6666
$ws.Range("M7").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 805
$ws.Range("I12").Value = 805
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 805
$ws.Range("L12").Value = 0
$ws.Range("H35").Value = 4999
$ws.Range("I35").Value = 4999
$ws.Range("K35").Value = 4999
$ws.Range("M35").Value = -4705
$ws.Range("H62").Value = 11113271
$ws.Range("I62").Value = 55555556
$ws.Range("J62").Value = 2700
$ws.Range("K62").Value = 55555556
$ws.Range("L62").Value = 2700
$ws.Range("M62").Value = -55554932
$ws.Range("N62").Value = -3948
$ws.Range("H65").Value = 11113271
$ws.Range("I65").Value = 55555556
$ws.Range("J65").Value = 2700
$ws.Range("K65").Value = 277777780
$ws.Range("L65").Value = 13500
$ws.Range("M65").Value = -277774660
$ws.Range("N65").Value = -19740
$ws.Range("H94").Value = 1657
$ws.Range("J94").Value = 1528.4
$ws.Range("L94").Value = 1528.4
$ws.Range("N94").Value = -2430.4
$ws.Range("H132").Value = 4800.4
$ws.Range("I132").Value = 4800.4
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 14401.2
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -11871.2
$ws.Range("M12").Value = -635
$ws.Range("N12").Value = ""
$ws.Range("N132").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 583.3333
$ws.Range("J68").Value = 583.3333
$ws.Range("L68").Value = 1749.9999
$ws.Range("N68").Value = -3371.9999
$ws.Range("H71").Value = 583.3333
$ws.Range("J71").Value = 583.3333
$ws.Range("L71").Value = 5249.9997
$ws.Range("N71").Value = -13361.9997
$ws.Range("H113").Value = 675.619
$ws.Range("I113").Value = 637.9167
$ws.Range("J113").Value = 725.8889
$ws.Range("K113").Value = 1913.7501
$ws.Range("L113").Value = 2177.6667
$ws.Range("M113").Value = 256.2499
$ws.Range("N113").Value = -6517.6667
$ws.Range("H131").Value = 875.76
$ws.Range("I131").Value = 740
$ws.Range("J131").Value = 881.4167
$ws.Range("K131").Value = 2220
$ws.Range("L131").Value = 2644.2501
$ws.Range("M131").Value = 2820
$ws.Range("N131").Value = -12724.2501

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 2669
$ws.Range("I29").Value = 2003.5
$ws.Range("J29").Value = 4000
$ws.Range("K29").Value = 2003.5
$ws.Range("L29").Value = 4000
$ws.Range("M29").Value = -1713.5
$ws.Range("N29").Value = -4580
$ws.Range("H122").Value = 3013.9412
$ws.Range("I122").Value = 2707.25
$ws.Range("J122").Value = 3750
$ws.Range("K122").Value = 8121.75
$ws.Range("L122").Value = 11250
$ws.Range("M122").Value = -5671.75
$ws.Range("N122").Value = -16150

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 61040.824
$ws.Range("I40").Value = 168415.67
$ws.Range("J40").Value = 2472.7273
$ws.Range("K40").Value = 168415.67
$ws.Range("L40").Value = 2472.7273
$ws.Range("M40").Value = -168279.67
$ws.Range("N40").Value = -2744.7273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 30151.727
$ws.Range("I56").Value = 3285
$ws.Range("J56").Value = 32838.4
$ws.Range("K56").Value = 3285
$ws.Range("L56").Value = 32838.4
$ws.Range("N56").Value = -34266.4
$ws.Range("H132").Value = 7751.143
$ws.Range("I132").Value = 8209.666999999999
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 24629.001
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -22099.001
$ws.Range("N132").Value = -20060
$ws.Range("M56").Value = -2571
